# "finish propensity to buy"
# Slide 2 holds the workflow diagram. Three plain-rectangle textboxes
# ("Propensity-to-buy", "Recommendation system", "Market basket analysis")
# are replaced by hand-drawn ("sketchy") textboxes matching the style
# already used for the other finished workflow items. The now-complete
# "Propensity-to-buy" item is renamed/re-styled (italic + accent2 fill) to
# match the other fully-finished items, while the other two keep the
# "in-progress" sketchy look. A handful of neighboring labels get nudged
# and italicized too.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function EmuToPt($emu) {
    return $emu / $EMU_PER_PT
}

# ---------------------------------------------------------------------
# 1. Drop the three "not started" plain textboxes - they get replaced by
#    sketchy versions below.
# ---------------------------------------------------------------------
$s.Shapes.Item("TextBox 5").Delete()   # "Propensity-to-buy"
$s.Shapes.Item("TextBox 6").Delete()   # "Recommendation system"
$s.Shapes.Item("TextBox 11").Delete()  # "Market basket analysis"

# ---------------------------------------------------------------------
# 2. Nudge a few existing shapes and make their labels italic (work in
#    progress -> highlighted).
# ---------------------------------------------------------------------

# "Comment sentiment prediction" (the upper one) shifts slightly.
$commentSentiment = $s.Shapes.Item("TextBox 31")
$commentSentiment.Left = EmuToPt 7302152
$commentSentiment.Top = EmuToPt 2694021

# "Product offering"
$productOffering = $s.Shapes.Item("TextBox 24")
$productOffering.Left = EmuToPt 1332720
$productOffering.Top = EmuToPt 1303274
$productOffering.TextFrame.TextRange.Font.Italic = -1

# "Personalized messaging"
$personalizedMessaging = $s.Shapes.Item("TextBox 36")
$personalizedMessaging.Left = EmuToPt 8344680
$personalizedMessaging.Top = EmuToPt 1283990
$personalizedMessaging.TextFrame.TextRange.ParagraphFormat.Alignment = 3 # ppAlignRight
$personalizedMessaging.TextFrame.TextRange.Font.Italic = -1

# "Preventing churn"
$preventingChurn = $s.Shapes.Item("TextBox 38")
$preventingChurn.Left = EmuToPt 8002209
$preventingChurn.Top = EmuToPt 5842847
$preventingChurn.TextFrame.TextRange.Font.Italic = -1

# "Profit enhancement"
$profitEnhancement = $s.Shapes.Item("TextBox 40")
$profitEnhancement.Left = EmuToPt 1332720
$profitEnhancement.Top = EmuToPt 5842847
$profitEnhancement.TextFrame.TextRange.Font.Italic = -1

# ---------------------------------------------------------------------
# 3. "Customer segmentation" moves from the light tx2 fill to accent2 -
#    this is the shared look for fully-finished items, and is also the
#    shape we borrow (duplicate) to create the new "Propensity-to-buy"
#    sketchy box below.
# ---------------------------------------------------------------------
$customerSeg = $s.Shapes.Item("TextBox 9")
$customerSeg.Fill.ForeColor.SchemeColor = "accent2"

# ---------------------------------------------------------------------
# 4. Build the new sketchy textboxes.
#    - "Propensity-to buy (PTB) model": finished -> duplicate the
#      (now accent2 + italic) "Customer segmentation" box.
#    - "Market basket analysis" / "Recommendation system": still in
#      progress -> duplicate the "Dynamic pricing" box (tx2 fill, not
#      italic, centered) which already carries the matching sketchy
#      style.
# ---------------------------------------------------------------------

# -- Propensity-to buy (PTB) model --
$dupRange = $customerSeg.Duplicate()
$ptb = $dupRange.Item(1)
$ptb.Name = "TextBox 18"
$ptb.Left = EmuToPt 2204541
$ptb.Top = EmuToPt 3180876
$ptb.Width = EmuToPt 3075709
$ptb.Height = EmuToPt 369332
$ptb.TextFrame.TextRange.Text = "Propensity-to buy (PTB) model"
$ptb.TextFrame.TextRange.Font.Italic = -1
$ptb.TextFrame.TextRange.ParagraphFormat.Alignment = 2 # ppAlignCenter

# -- Market basket analysis --
$dynamicPricing = $s.Shapes.Item("TextBox 41")
$dupRange = $dynamicPricing.Duplicate()
$mba = $dupRange.Item(1)
$mba.Name = "TextBox 19"
$mba.Left = EmuToPt 1791877
$mba.Top = EmuToPt 4053603
$mba.Width = EmuToPt 2424279
$mba.Height = EmuToPt 369332
$mba.TextFrame.TextRange.Text = "Market basket analysis"
$mba.TextFrame.TextRange.ParagraphFormat.Alignment = 2 # ppAlignCenter

# -- Recommendation system --
$dupRange = $dynamicPricing.Duplicate()
$recSys = $dupRange.Item(1)
$recSys.Name = "TextBox 20"
$recSys.Left = EmuToPt 1791877
$recSys.Top = EmuToPt 2324689
$recSys.Width = EmuToPt 2704402
$recSys.Height = EmuToPt 369332
$recSys.TextFrame.TextRange.Text = "Recommendation system"
$recSys.TextFrame.TextRange.ParagraphFormat.Alignment = 2 # ppAlignCenter
